$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 47: this row already existed with only Date/Travail effectué; fill
# in the Type and Temps [h] columns that were left blank. ---
$ws.Range("B47").Value = "Rédaction"
$ws.Range("C47").Value = 3

# --- Row 48: brand-new journal entry ---
$ws.Range("A48").Value = 45076
$ws.Range("B48").Value = "Implémentation"
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = "Frontend: Interface config module"

# The old totals row (row 49, "=SUM(C2:C48)") needs to move down to make
# room for the new blank rows added to the table; clear it now so it
# doesn't leak into the new placeholder rows below.
$ws.Range("C49").ClearContents()

# --- Give the new date cell (A48) the same style as the rest of column A,
# then stamp that same format down through the blank placeholder rows
# (49-60) that the table grows into. ---
$ws.Range("A47").Copy()
$ws.Range("A48:A60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column C (Temps [h]) placeholder styling: round-trip the number
# format on C48 so the cell carries an explicit (General) number format,
# then copy that style down through the rest of the blank rows. ---
$ws.Range("C48").NumberFormat = "0.00"
$ws.Range("C48").NumberFormat = "General"
$ws.Range("C48").Copy()
$ws.Range("C48:C60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 61 is left completely empty (no data, no explicit formatting) - the
# table simply grows over it without any cell content.

# --- Row 62: new totals row for the table ---
$ws.Range("C62").Formula = "=SUM(C2:C61)"

# --- Resize/extend the table to cover the new data + totals row, which
# also grows the AutoFilter range accordingly. ---
$lo.Resize($ws.Range("A1:E62"))

# --- Move the active selection the same way the author's workbook shows
# it after the edit. ---
$ws.Range("D50").Select()

Write-Host "edit applied"
